# Commit: Rename sheet 'Data' to 'Data table' (Close #151)
#
# The author renamed the "Data" worksheet to "Data table" in Excel. Doing
# this interactively also left that sheet as the active/selected tab (it
# was clicked on to trigger the rename / inspect the result), while the
# sheet that used to be selected ("Table") kept its own last-used cell
# selection, which had moved on (B17 -> B89) before the save.

$wb = $excel.ActiveWorkbook

# Preserve/update the selection on the previously active sheet ("Table")
# before switching away from it, so its stored <selection> reflects where
# the cursor ended up (B17 -> B89).
$wsTable = $wb.Worksheets.Item("Table")
[void]$wsTable.Range("B89").Select()

# Rename "Data" -> "Data table".
$wsData = $wb.Worksheets.Item("Data")
$wsData.Name = "Data table"

# Make the renamed sheet the active tab (it keeps its own selection at D1).
[void]$wsData.Activate()
[void]$wsData.Range("D1").Select()
